$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Exclude internal from formatted": two internal/duplicate entries
# (originally on rows 108 and 109 - an "Adm support (USAID/OFDA)" style
# record and the "Support to: (i) treatment of Moderate Acute
# Malnutrition (MAM) ..." record) are not meant to ship in the
# formatted output. Remove those two data rows entirely so every row
# below shifts up by two; the now-unused "MAM" shared string is dropped
# from the workbook automatically since nothing references it anymore.
$ws.Rows("108:109").Delete()
